$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H), mirroring the style of the existing
# header cells (bold, centered, thin border) by copying G1's format.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New data values for the "Save" column.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
